# Generate Report for Handback
#
# The 39a73dd1-26b3-40d2-bbc0-a18f30b162eb handback row (row 6) on both the
# "zh-cn" and "de-de" language sheets gets its localization-report columns
# filled in (Latest Target File / Latest Handback File / Latest Handback
# DateTime / Error Detail), because the handed-back file turned out not to
# be based on the latest handoff version. The "Error Detail" column (P) is
# also widened so the message is readable.

$wb = $excel.ActiveWorkbook

$targetMdName = "39a73dd1-26b3-40d2-bbc0-a18f30b162eb.md"
$targetMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bda951e2730eb638ac91cf3219c793f4defe1348/e2e/39a73dd1-26b3-40d2-bbc0-a18f30b162eb.md"
$errorDetail  = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aff22a4bb7b1a6361dbc618bcce03b82fb8d2aa3/e2e/39a73dd1-26b3-40d2-bbc0-a18f30b162eb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bda951e2730eb638ac91cf3219c793f4defe1348/e2e/39a73dd1-26b3-40d2-bbc0-a18f30b162eb.md."

# Cornflower blue (FF6495ED), same color already used by the workbook's
# "HyperLink" cell style, expressed as the BGR integer the Font.Color
# COM property expects.
$hyperlinkColor = 15570276

# NOTE: this runtime's PowerShell-style function calls don't bind
# `-Name value` style named arguments reliably, so the helper below takes
# plain positional parameters.
function Set-HandbackReportRow($Worksheet, [string]$HandbackXlf, [string]$HandbackDateTime) {
    # I6: Latest Target File -> the handoff markdown, now a hyperlink.
    $Worksheet.Range("I6").Value = $targetMdName
    $Worksheet.Range("I6").Font.Underline = 2
    $Worksheet.Range("I6").Font.Color = $hyperlinkColor
    $Worksheet.Hyperlinks.Add($Worksheet.Range("I6"), $targetMdUrl, $null, $null, $targetMdName) | Out-Null

    # J6: Latest Handback File -> the xlf that was handed back.
    $Worksheet.Range("J6").Value = $HandbackXlf

    # K6: Latest Handback DateTime -> when that handback happened.
    $Worksheet.Range("K6").Value = $HandbackDateTime

    # P6: Error Detail -> explains the handback wasn't against the latest version.
    $Worksheet.Range("P6").Value = $errorDetail
}

# The stored OOXML column width (40) and the COM `ColumnWidth` the host
# exposes differ by the fixed ~5/6 character "padding" this engine (like
# Excel itself) adds when round-tripping character-width <-> pixel-width;
# column A already stores width=40 and reads back as ColumnWidth 39.17, so
# we use that same value here to land on an OOXML width of exactly 40.
$errorDetailColumnWidth = 39.166666666666664

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackReportRow $wsZhCn "39a73dd1-26b3-40d2-bbc0-a18f30b162eb.f9938f3ea1c8bcc47da59dbbec377b434bead967.zh-cn.xlf" "2016-09-05 18:48:43"
$wsZhCn.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackReportRow $wsDeDe "39a73dd1-26b3-40d2-bbc0-a18f30b162eb.f9938f3ea1c8bcc47da59dbbec377b434bead967.de-de.xlf" "2016-09-05 18:48:51"
$wsDeDe.Columns.Item(16).ColumnWidth = $errorDetailColumnWidth
